# Auto-generated edit script applying the diff to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("N2").Value = 1.3
$ws.Range("P2").Value = 1.3
$ws.Range("S2").Value = 1.37

# Row 3
$ws.Range("N3").Value = 1.27
$ws.Range("P3").Value = 1.27
$ws.Range("Q3").Value = 1.02
$ws.Range("S3").Value = 1.39

# Row 5
$ws.Range("F5").Value = 1.62
$ws.Range("I5").Value = 6
$ws.Range("L5").Value = 1.01
$ws.Range("N5").Value = 5.9
$ws.Range("O5").Value = 1.15
$ws.Range("P5").Value = 2.76
$ws.Range("Q5").Value = 1.46
$ws.Range("R5").Value = 1.73
$ws.Range("S5").Value = 2.16
$ws.Range("T5").Value = 1.56
$ws.Range("V5").Value = 1.22
$ws.Range("X5").Value = 34
$ws.Range("AA5").Value = 150
$ws.Range("AC5").Value = 13
$ws.Range("AN5").Value = 6.4

# Row 6
$ws.Range("F6").Value = 2.74
$ws.Range("G6").Value = 3.05
$ws.Range("L6").Value = 1.34
$ws.Range("N6").Value = 3.7
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 1.45
$ws.Range("S6").Value = 2.56
$ws.Range("T6").Value = 1.63
$ws.Range("W6").Value = 1.48
$ws.Range("X6").Value = 19
$ws.Range("Y6").Value = 13.5
$ws.Range("Z6").Value = 19.5
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 14.5
$ws.Range("AC6").Value = 9
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 28
$ws.Range("AF6").Value = 22
$ws.Range("AG6").Value = 14
$ws.Range("AH6").Value = 16.5
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 24
$ws.Range("AO6").Value = 19

# Row 7
$ws.Range("F7").Value = 3.35
$ws.Range("H7").Value = 2.1
$ws.Range("I7").Value = 2.2
$ws.Range("J7").Value = 3.7
$ws.Range("K7").Value = 4.4
$ws.Range("N7").Value = 4.8
$ws.Range("O7").Value = 1.17
$ws.Range("P7").Value = 2.32
$ws.Range("R7").Value = 1.55
$ws.Range("S7").Value = 2.06
$ws.Range("T7").Value = 1.52
$ws.Range("U7").Value = 2.66
$ws.Range("V7").Value = 1.83
$ws.Range("Y7").Value = 1000
$ws.Range("AC7").Value = 10.5
$ws.Range("AH7").Value = 15
$ws.Range("AO7").Value = 9.4

# Row 8
$ws.Range("H8").Value = 1.51
$ws.Range("I8").Value = 1.53
$ws.Range("J8").Value = 4.7
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 1.28
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 1.19
$ws.Range("Q8").Value = 1.59
$ws.Range("R8").Value = 1.6
$ws.Range("S8").Value = 2.48
$ws.Range("T8").Value = 1.76
$ws.Range("U8").Value = 2.12
$ws.Range("V8").Value = 2.84
$ws.Range("X8").Value = 24
$ws.Range("Y8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 38
$ws.Range("AD8").Value = 980
$ws.Range("AG8").Value = 26
$ws.Range("AI8").Value = 980
$ws.Range("AK8").Value = 1000
$ws.Range("AN8").Value = 95
$ws.Range("AO8").Value = 1000
